# live_trading_results.xlsx update
# Trade #137 (row 166 on "All Trades" / row 25 on "HighProbConvergence")
# closed at 2026-02-18 00:41:33, plus three new OPEN trades appended
# (momentum #194, MarketMaking #195, EMAArbitrage #196) and the
# Summary / Strategy Status roll-up numbers refreshed accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper: write a date-like text value (e.g. "2026-02-18") without
# letting the host auto-convert it into a real date serial/style.
# A leading apostrophe forces text, then the style is reset to
# "Normal" so no number-format style sticks to the cell.
# ---------------------------------------------------------------
function Set-TextCell($cell, [string]$text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# =================================================================
# Sheet "Summary"
# =================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.12   # Current Capital
$summary.Range("B4").Value = 0.23      # Total P&L $
$summary.Range("B6").Value = 165       # Total Trades
$summary.Range("B8").Value = 59        # Losing Trades
$summary.Range("B9").Value = 45.45     # Win Rate %

# =================================================================
# Sheet "Strategy Status" -- row 3 = HighProbConvergence
# =================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.34
$status.Range("D3").Value = 24
$status.Range("E3").Value = 0.35
$status.Range("F3").Value = 0.34
$status.Range("G3").Value = 58.33

# =================================================================
# Sheet "All Trades" -- close trade #165 (row 166)
# =================================================================
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(166, 7).Value = 0.82                 # G166 Exit Price
$allTrades.Cells.Item(166, 8).Value = "CLOSED"              # H166 Status
$allTrades.Cells.Item(166, 9).Value = -3.5294                # I166 P&L %
$allTrades.Cells.Item(166, 10).Value = -0.03                 # J166 P&L $
$allTrades.Cells.Item(166, 11).Value = 100.34                # K166 Capital After
$allTrades.Cells.Item(166, 12).Value = "early_exit"          # L166 Entry Slippage(label col here is text)
$allTrades.Cells.Item(166, 13).Value = 0.19                  # M166 Exit Slippage

# =================================================================
# Sheet "HighProbConvergence" -- mirror of the same trade (row 25)
# =================================================================
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(25, 7).Value = 0.82                  # G25 Exit Price
$hpc.Cells.Item(25, 8).Value = "CLOSED"               # H25 Status
$hpc.Cells.Item(25, 9).Value = -3.5294                 # I25 P&L %
$hpc.Cells.Item(25, 10).Value = -0.03                  # J25 P&L $
$hpc.Cells.Item(25, 11).Value = 100.34                 # K25 Capital After
$hpc.Cells.Item(25, 16).Value = "early_exit"           # P25 Exit Reason
$hpc.Cells.Item(25, 17).Value = 0.19                   # Q25 Duration (min)

# =================================================================
# New trade rows (OPEN) appended to "All Trades" and to each
# strategy's own sheet.
# =================================================================

# ---- All Trades row 195 : Trade #194, momentum, DOWN ----
$r = $allTrades.Cells.Item(195, 1)
$r.Value = 194
Set-TextCell $allTrades.Cells.Item(195, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(195, 3) "00:41:26"
$allTrades.Cells.Item(195, 4).Value = "momentum"
$allTrades.Cells.Item(195, 5).Value = "DOWN"
$allTrades.Cells.Item(195, 6).Value = 0.85
$allTrades.Cells.Item(195, 8).Value = "OPEN"
$allTrades.Cells.Item(195, 9).Value = 0
$allTrades.Cells.Item(195, 10).Value = 0
$allTrades.Cells.Item(195, 11).Value = 99.19712996249174
$allTrades.Cells.Item(195, 13).Value = 0
$allTrades.Cells.Item(195, 14).Value = 0
$allTrades.Cells.Item(195, 15).Value = 0
$allTrades.Cells.Item(195, 16).Value = 0.9
$allTrades.Cells.Item(195, 17).Value = "Downward momentum: -45.109% over 10 samples"

# ---- All Trades row 196 : Trade #195, MarketMaking, UP ----
$allTrades.Cells.Item(196, 1).Value = 195
Set-TextCell $allTrades.Cells.Item(196, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(196, 3) "00:41:27"
$allTrades.Cells.Item(196, 4).Value = "MarketMaking"
$allTrades.Cells.Item(196, 5).Value = "UP"
$allTrades.Cells.Item(196, 6).Value = 0.15
$allTrades.Cells.Item(196, 8).Value = "OPEN"
$allTrades.Cells.Item(196, 9).Value = 0
$allTrades.Cells.Item(196, 10).Value = 0
$allTrades.Cells.Item(196, 11).Value = 99.30858346467944
$allTrades.Cells.Item(196, 13).Value = 0
$allTrades.Cells.Item(196, 14).Value = 0
$allTrades.Cells.Item(196, 15).Value = 0
$allTrades.Cells.Item(196, 16).Value = 0.6
$allTrades.Cells.Item(196, 17).Value = "Normal spread capture: 198 bps"

# ---- All Trades row 197 : Trade #196, EMAArbitrage, DOWN ----
$allTrades.Cells.Item(197, 1).Value = 196
Set-TextCell $allTrades.Cells.Item(197, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(197, 3) "00:41:27"
$allTrades.Cells.Item(197, 4).Value = "EMAArbitrage"
$allTrades.Cells.Item(197, 5).Value = "DOWN"
$allTrades.Cells.Item(197, 6).Value = 0.84
$allTrades.Cells.Item(197, 8).Value = "OPEN"
$allTrades.Cells.Item(197, 9).Value = 0
$allTrades.Cells.Item(197, 10).Value = 0
$allTrades.Cells.Item(197, 11).Value = 100.270616878256
$allTrades.Cells.Item(197, 13).Value = 0
$allTrades.Cells.Item(197, 14).Value = 0
$allTrades.Cells.Item(197, 15).Value = 0
$allTrades.Cells.Item(197, 16).Value = 0.9
$allTrades.Cells.Item(197, 17).Value = "EMA:down, RSI:50.0, ROC:-45.11% | 2/3 DOWN"

# ---- "momentum" sheet row 53 : Trade #194 ----
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(53, 1).Value = 194
Set-TextCell $momentum.Cells.Item(53, 2) "2026-02-18"
Set-TextCell $momentum.Cells.Item(53, 3) "00:41:26"
$momentum.Cells.Item(53, 4).Value = "momentum"
$momentum.Cells.Item(53, 5).Value = "DOWN"
$momentum.Cells.Item(53, 6).Value = 0.85
$momentum.Cells.Item(53, 8).Value = "OPEN"
$momentum.Cells.Item(53, 9).Value = 0
$momentum.Cells.Item(53, 10).Value = 0
$momentum.Cells.Item(53, 11).Value = 99.19712996249174
$momentum.Cells.Item(53, 12).Value = 0
$momentum.Cells.Item(53, 13).Value = 0
$momentum.Cells.Item(53, 14).Value = 0.9
$momentum.Cells.Item(53, 15).Value = "Downward momentum: -45.109% over 10 samples"
$momentum.Cells.Item(53, 17).Value = 0

# ---- "MarketMaking" sheet row 83 : Trade #195 ----
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(83, 1).Value = 195
Set-TextCell $marketMaking.Cells.Item(83, 2) "2026-02-18"
Set-TextCell $marketMaking.Cells.Item(83, 3) "00:41:27"
$marketMaking.Cells.Item(83, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(83, 5).Value = "UP"
$marketMaking.Cells.Item(83, 6).Value = 0.15
$marketMaking.Cells.Item(83, 8).Value = "OPEN"
$marketMaking.Cells.Item(83, 9).Value = 0
$marketMaking.Cells.Item(83, 10).Value = 0
$marketMaking.Cells.Item(83, 11).Value = 99.30858346467944
$marketMaking.Cells.Item(83, 12).Value = 0
$marketMaking.Cells.Item(83, 13).Value = 0
$marketMaking.Cells.Item(83, 14).Value = 0.6
$marketMaking.Cells.Item(83, 15).Value = "Normal spread capture: 198 bps"
$marketMaking.Cells.Item(83, 17).Value = 0

# ---- "EMAArbitrage" sheet row 9 : Trade #196 ----
$emaArb = $wb.Worksheets.Item("EMAArbitrage")
$emaArb.Cells.Item(9, 1).Value = 196
Set-TextCell $emaArb.Cells.Item(9, 2) "2026-02-18"
Set-TextCell $emaArb.Cells.Item(9, 3) "00:41:27"
$emaArb.Cells.Item(9, 4).Value = "EMAArbitrage"
$emaArb.Cells.Item(9, 5).Value = "DOWN"
$emaArb.Cells.Item(9, 6).Value = 0.84
$emaArb.Cells.Item(9, 8).Value = "OPEN"
$emaArb.Cells.Item(9, 9).Value = 0
$emaArb.Cells.Item(9, 10).Value = 0
$emaArb.Cells.Item(9, 11).Value = 100.270616878256
$emaArb.Cells.Item(9, 12).Value = 0
$emaArb.Cells.Item(9, 13).Value = 0
$emaArb.Cells.Item(9, 14).Value = 0.9
$emaArb.Cells.Item(9, 15).Value = "EMA:down, RSI:50.0, ROC:-45.11% | 2/3 DOWN"
$emaArb.Cells.Item(9, 17).Value = 0
